$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.279.22"
$ws.Range("E2").Value = "  +1.23%  "
$ws.Range("D3").Value = "1.610.12"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.99"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("E7").Value = "  +0.53%  "
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.42"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0815"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.42%  "
$ws.Range("D12").Value = "1.835.01"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("D13").Value = "1.609.74"
$ws.Range("E13").Value = "  +1.04%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "26.271.21"
$ws.Range("E16").Value = "  +1.12%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.07"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.09%  "
$ws.Range("E18").Value = "  +0.89%  "
$ws.Range("E19").Value = "  -0.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "201.26"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.32"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.91"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.62%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.24"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.24"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("E29").Value = "  +2.39%  "
$ws.Range("E30").Value = "  +6.06%  "
$ws.Range("E31").Value = "  +0.24%  "
$ws.Range("E32").Value = "  +2.63%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  +1.21%  "
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "1.157.20"
$ws.Range("E36").Value = "  +2.72%  "
$ws.Range("E37").Value = "  +0.96%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +1.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.789"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("E41").Value = "  +1.50%  "
$ws.Range("E42").Value = "  +4.24%  "
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("D44").Value = "1.746.06"
$ws.Range("E44").Value = "  +0.60%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.84"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.25%  "
$ws.Range("E46").Value = "  +14.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.52"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("E49").Value = "  +0.85%  "
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("E51").Value = "  -0.41%  "
